$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineShapeXmlName($shapeRange, $oldName, $newName) {
    $xml = $shapeRange.WordOpenXML
    $search = 'name="' + $oldName + '"'
    $replace = 'name="' + $newName + '"'
    $updated = $xml.Replace($search, $replace)
    $shapeRange.WordOpenXML = $updated
}

# Footer (Primary / default) -> footer2.xml, docPr id="2", PearsonLogo: image2.png -> image1.png
$ftrPrimary = $sec.Footers.Item(1)
$shp = $ftrPrimary.Range.InlineShapes.Item(1)
Rename-InlineShapeXmlName $shp.Range "image2.png" "image1.png"

# Footer (First Page) -> footer1.xml, docPr id="3", PearsonLogo: image2.png -> image1.png
$ftrFirst = $sec.Footers.Item(2)
$shp = $ftrFirst.Range.InlineShapes.Item(1)
Rename-InlineShapeXmlName $shp.Range "image2.png" "image1.png"

# Header (First Page) -> header1.xml, docPr id="1", BTec_Logo-Orange: image1.jpg -> image2.jpg
$hdrFirst = $sec.Headers.Item(2)
$shp = $hdrFirst.Range.InlineShapes.Item(1)
Rename-InlineShapeXmlName $shp.Range "image1.jpg" "image2.jpg"
